$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update last-updated timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "30 Oct 2025, 12:22 PM"

# --- "1 Month Performance" sheet: update stock names / % change values ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Range("C4").Value = 78.68340000000001
$wsPerf.Range("C5").Value = 66.4308
$wsPerf.Range("C6").Value = 63.1026
$wsPerf.Range("C7").Value = 57.0641
$wsPerf.Range("C9").Value = 53.7803
$wsPerf.Range("C10").Value = 45.435
$wsPerf.Range("C11").Value = 42.5017
$wsPerf.Range("C14").Value = 38.8863
$wsPerf.Range("C15").Value = 38.6706
$wsPerf.Range("C16").Value = 37.3091
$wsPerf.Range("C17").Value = 36.5209
$wsPerf.Range("C19").Value = 34.5742
$wsPerf.Range("C20").Value = 34.3103
$wsPerf.Range("C22").Value = 33.3711
$wsPerf.Range("C23").Value = 32.6637
$wsPerf.Range("C24").Value = 30.5748
$wsPerf.Range("C26").Value = 29.8891
$wsPerf.Range("C28").Value = 29.4535
$wsPerf.Range("B29").Value = "ARFIN"
$wsPerf.Range("C29").Value = 28.6249
$wsPerf.Range("B30").Value = "TARACHAND"
$wsPerf.Range("C30").Value = 28.4897
$wsPerf.Range("C31").Value = 27.1044
$wsPerf.Range("C32").Value = 26.5533
$wsPerf.Range("C33").Value = 26.2677
$wsPerf.Range("B34").Value = "EMKAY"
$wsPerf.Range("C34").Value = 25.7923
$wsPerf.Range("B35").Value = "SAGILITY"
$wsPerf.Range("C35").Value = 25.6204
$wsPerf.Range("C36").Value = 25.2799
$wsPerf.Range("C37").Value = 25.2459
$wsPerf.Range("B38").Value = "MARINE"
$wsPerf.Range("C38").Value = 25.2041
$wsPerf.Range("B39").Value = "AVALON"
$wsPerf.Range("C39").Value = 24.8392
$wsPerf.Range("C41").Value = 24.5059
$wsPerf.Range("C42").Value = 24.3737
$wsPerf.Range("C43").Value = 23.8751
$wsPerf.Range("B46").Value = "LORDSCHLO"
$wsPerf.Range("C46").Value = 23.6702
$wsPerf.Range("B47").Value = "CARTRADE"
$wsPerf.Range("C47").Value = 23.5413
$wsPerf.Range("B48").Value = "TATVA"
$wsPerf.Range("C48").Value = 22.9632
$wsPerf.Range("B49").Value = "INDIANB"
$wsPerf.Range("C49").Value = 22.6452
$wsPerf.Range("C50").Value = 22.3662
$wsPerf.Range("C51").Value = 22.2685
$wsPerf.Range("C52").Value = 22.0372
$wsPerf.Range("B54").Value = "GUJTHEM"
$wsPerf.Range("C54").Value = 21.718
$wsPerf.Range("B55").Value = "IIFL"
$wsPerf.Range("C55").Value = 21.6364
$wsPerf.Range("C60").Value = 20.353
$wsPerf.Range("B61").Value = "FEDERALBNK"
$wsPerf.Range("C61").Value = 20.1975
$wsPerf.Range("B62").Value = "HINDCOPPER"
$wsPerf.Range("C62").Value = 20.1762
$wsPerf.Range("B63").Value = "GRMOVER"
$wsPerf.Range("C63").Value = 19.9714
$wsPerf.Range("B64").Value = "BHARATWIRE"
$wsPerf.Range("C64").Value = 19.8904
$wsPerf.Range("B65").Value = "BHAGERIA"
$wsPerf.Range("C65").Value = 19.5958
$wsPerf.Range("B66").Value = "SHRIRAMFIN"
$wsPerf.Range("C66").Value = 19.5599
$wsPerf.Range("C67").Value = 19.3568
$wsPerf.Range("B68").Value = "MCX"
$wsPerf.Range("C68").Value = 19.1289
$wsPerf.Range("B69").Value = "CEATLTD"
$wsPerf.Range("C69").Value = 19.069
$wsPerf.Range("B72").Value = "REPRO"
$wsPerf.Range("C72").Value = 18.6348
$wsPerf.Range("B73").Value = "WHEELS"
$wsPerf.Range("C73").Value = 18.4184

# --- "distance from Dma50" sheet: update distance values ---
$wsDma = $wb.Worksheets.Item("distance from Dma50")
$wsDma.Range("C2").Value = 9.528499999999999
$wsDma.Range("C3").Value = 7.3726
$wsDma.Range("C4").Value = 6.3763
$wsDma.Range("C5").Value = 5.3602
$wsDma.Range("C6").Value = 5.2639
$wsDma.Range("C7").Value = 5.0117
$wsDma.Range("C8").Value = 4.4151
$wsDma.Range("C9").Value = 4.3828
$wsDma.Range("C10").Value = 3.8381
$wsDma.Range("C11").Value = 3.6543
$wsDma.Range("C12").Value = 3.3799
$wsDma.Range("C13").Value = 3.3769
$wsDma.Range("C15").Value = 3.0332
$wsDma.Range("C16").Value = 2.9509
$wsDma.Range("C17").Value = 2.8154
$wsDma.Range("C18").Value = 2.7877
$wsDma.Range("C19").Value = 2.7105
$wsDma.Range("C20").Value = 2.3479
$wsDma.Range("C21").Value = 2.3102
$wsDma.Range("C22").Value = 1.4094
$wsDma.Range("C23").Value = 1.3991
$wsDma.Range("C24").Value = 1.2713
$wsDma.Range("C25").Value = 1.0724
$wsDma.Range("C26").Value = 0.9967
$wsDma.Range("C27").Value = 0.8801
$wsDma.Range("C28").Value = 0.4843
$wsDma.Range("C29").Value = 0.3324
$wsDma.Range("C30").Value = -2.0251
